# Existing API updated with latest payload
# Update marketplace rows with the latest payload values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: QuantumSoft -> Digicap
$ws.Range("B2").Value = "Digicap"
$ws.Range("H2").Value = "digicap.aidtaas.ai"

# Row 3: FashionHub -> Amazon
$ws.Range("A3").Value = "Error: 201"
$ws.Range("B3").Value = "Amazon"
$ws.Range("H3").Value = "amazon.aidtaas.ai"

# Row 4: FinSecure -> HP
$ws.Range("A4").Value = "Error: 201"
$ws.Range("B4").Value = "HP"
$ws.Range("H4").Value = "hp.aidtaas.ai"
